# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) on several rows across multiple sheets, reflecting refreshed
# market-board data from the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 694.96
$ws.Range("I19").Value2 = 667.8823
$ws.Range("J19").Value2 = 752.5
$ws.Range("K19").Value2 = 667.8823
$ws.Range("L19").Value2 = 752.5
$ws.Range("M19").Value2 = -492.8823
$ws.Range("N19").Value2 = -1102.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value2 = 2254.3
$ws.Range("J88").Value2 = 2480.375
$ws.Range("L88").Value2 = 2480.375
$ws.Range("N88").Value2 = -3292.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value2 = 2254.3
$ws.Range("J91").Value2 = 2480.375
$ws.Range("L91").Value2 = 2480.375
$ws.Range("N91").Value2 = -5288.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value2 = 83334680
$ws.Range("I111").Value2 = 512
$ws.Range("J111").Value2 = 142859090
$ws.Range("K111").Value2 = 1536
$ws.Range("L111").Value2 = 428577270
$ws.Range("M111").Value2 = 1531
$ws.Range("N111").Value2 = -428583404

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value2 = 1438.8462
$ws.Range("I132").Value2 = 1548
$ws.Range("K132").Value2 = 4644
$ws.Range("M132").Value2 = -2114

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 976.8
$ws.Range("I137").Value2 = 904.02563
$ws.Range("J137").Value2 = 1234.8182
$ws.Range("K137").Value2 = 2712.07689
$ws.Range("L137").Value2 = 3704.4546
$ws.Range("M137").Value2 = -162.0768899999998
$ws.Range("N137").Value2 = -8804.454600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2685.14
$ws.Range("I32").Value2 = 2680.763
$ws.Range("J32").Value2 = 2826.6667
$ws.Range("K32").Value2 = 2680.763
$ws.Range("L32").Value2 = 2826.6667
$ws.Range("M32").Value2 = -2393.763
$ws.Range("N32").Value2 = -3400.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 973.5
$ws.Range("I61").Value2 = 973.5
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 973.5
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -761.5
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3309.2917
$ws.Range("I74").Value2 = 3473.4546
$ws.Range("K74").Value2 = 3473.4546
$ws.Range("M74").Value2 = -2599.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 3309.2917
$ws.Range("I77").Value2 = 3473.4546
$ws.Range("K77").Value2 = 17367.273
$ws.Range("M77").Value2 = -12999.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 973.5
$ws.Range("I136").Value2 = 973.5
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 2920.5
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -370.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 1074.2545
$ws.Range("I134").Value2 = 829.0213
$ws.Range("K134").Value2 = 2487.0639
$ws.Range("M134").Value2 = 47.9360999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2849.776
$ws.Range("I31").Value2 = 1975.7391
$ws.Range("J31").Value2 = 3424.1428
$ws.Range("K31").Value2 = 1975.7391
$ws.Range("L31").Value2 = 3424.1428
$ws.Range("M31").Value2 = -1680.7391
$ws.Range("N31").Value2 = -4014.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 2849.776
$ws.Range("I34").Value2 = 1975.7391
$ws.Range("J34").Value2 = 3424.1428
$ws.Range("K34").Value2 = 1975.7391
$ws.Range("L34").Value2 = 3424.1428
$ws.Range("M34").Value2 = -1773.7391
$ws.Range("N34").Value2 = -3828.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 1343.4492
$ws.Range("I58").Value2 = 1037.3889
$ws.Range("K58").Value2 = 1037.3889
$ws.Range("M58").Value2 = -834.3888999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value2 = 1821.6666
$ws.Range("I105").Value2 = 1986
$ws.Range("J105").Value2 = 1000
$ws.Range("K105").Value2 = 1986
$ws.Range("L105").Value2 = 1000
$ws.Range("M105").Value2 = -239
$ws.Range("N105").Value2 = -4494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 1700.8055
$ws.Range("I132").Value2 = 908.4815
$ws.Range("J132").Value2 = 4077.7778
$ws.Range("K132").Value2 = 2725.4445
$ws.Range("L132").Value2 = 12233.3334
$ws.Range("M132").Value2 = -195.4445000000001
$ws.Range("N132").Value2 = -17293.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 1327.746
$ws.Range("I134").Value2 = 1138.9474
$ws.Range("K134").Value2 = 3416.8422
$ws.Range("M134").Value2 = -881.8422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 1343.4492
$ws.Range("I136").Value2 = 1037.3889
$ws.Range("K136").Value2 = 3112.1667
$ws.Range("M136").Value2 = -562.1666999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value2 = 2682.7368
$ws.Range("I115").Value2 = 744.5
$ws.Range("J115").Value2 = 3199.6
$ws.Range("K115").Value2 = 2233.5
$ws.Range("L115").Value2 = 9598.799999999999
$ws.Range("M115").Value2 = -1058.5
$ws.Range("N115").Value2 = -11948.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value2 = 1565.8064
$ws.Range("J132").Value2 = 1489.0834
$ws.Range("L132").Value2 = 13401.7506
$ws.Range("N132").Value2 = -18461.7506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 1268.8889
$ws.Range("I107").Value2 = 2190
$ws.Range("J107").Value2 = 117.5
$ws.Range("K107").Value2 = 2190
$ws.Range("L107").Value2 = 117.5
$ws.Range("M107").Value2 = -270
$ws.Range("N107").Value2 = -3957.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 2222.7144
$ws.Range("I113").Value2 = 2093.1667
$ws.Range("K113").Value2 = 2093.1667
$ws.Range("M113").Value2 = 76.83329999999978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value2 = 10271.667
$ws.Range("J123").Value2 = 10271.667
$ws.Range("L123").Value2 = 10271.667
$ws.Range("N123").Value2 = -15171.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 2625
$ws.Range("I61").Value2 = 2500
$ws.Range("K61").Value2 = 2500
$ws.Range("M61").Value2 = -2298

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value2 = 2625
$ws.Range("I113").Value2 = 2500
$ws.Range("K113").Value2 = 2500
$ws.Range("M113").Value2 = -330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 4924.3413
$ws.Range("I132").Value2 = 4583.183
$ws.Range("J132").Value2 = 5854.773
$ws.Range("K132").Value2 = 13749.549
$ws.Range("L132").Value2 = 17564.319
$ws.Range("M132").Value2 = -11219.549
$ws.Range("N132").Value2 = -22624.319

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1355.8206
$ws.Range("I132").Value2 = 1355.8206
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 4067.4618
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -1537.4618
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 1335.5962
$ws.Range("I136").Value2 = 500.93878
$ws.Range("J136").Value2 = 14968.333
$ws.Range("K136").Value2 = 1502.81634
$ws.Range("L136").Value2 = 44904.999
$ws.Range("M136").Value2 = 1047.18366
$ws.Range("N136").Value2 = -50004.999
